# -----------------------------------------------------------------------
# CryCompanywiseStockReport_1.xlsx - stock-quantity correction pass.
#
# The report lists one row per (company, item): purchase rate (D), sale
# rate (E), closing qty (F) and closing value (G = D * F). Per-company
# "Sub Total:" rows sum G for that company, and the sheet-wide
# "Sub Total:"/"Grand Total:" rows (B942/B943) sum every company subtotal.
#
# This pass re-counts several items closing stock. For plain corrections
# we only touch Qty (F) and re-derive Value (G = D * F). A few pairs of
# rows for the same item (different batches/rates) had their stock
# entries swapped to the other row. Every subtotal/grand total affected
# by these changes is then refreshed to match.
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Quantity (F) / Value (G) corrections -------------------------------
$ws.Range("F27").Value2 = 29
$ws.Range("G27").Value2 = 1470.88
$ws.Range("F32").Value2 = 115
$ws.Range("G32").Value2 = 4123.9
$ws.Range("F36").Value2 = 232
$ws.Range("G36").Value2 = 6537.76
$ws.Range("F37").Value2 = 146
$ws.Range("G37").Value2 = 3889.44
$ws.Range("F48").Value2 = 272
$ws.Range("G48").Value2 = 52466.08
$ws.Range("F51").Value2 = 51
$ws.Range("G51").Value2 = 1801.32
$ws.Range("F55").Value2 = 273
$ws.Range("G55").Value2 = 15315.3
$ws.Range("F59").Value2 = 218
$ws.Range("G59").Value2 = 20391.72
$ws.Range("F60").Value2 = 150
$ws.Range("G60").Value2 = 8850
$ws.Range("F63").Value2 = 165
$ws.Range("G63").Value2 = 9200.4
$ws.Range("F66").Value2 = 142
$ws.Range("G66").Value2 = 11066.06
$ws.Range("F68").Value2 = 22
$ws.Range("G68").Value2 = 411.62
$ws.Range("F69").Value2 = 284
$ws.Range("G69").Value2 = 74047.32000000001
$ws.Range("F99").Value2 = 7
$ws.Range("G99").Value2 = 1569.89
$ws.Range("F131").Value2 = 68
$ws.Range("G131").Value2 = 4331.6
$ws.Range("F141").Value2 = 23
$ws.Range("G141").Value2 = 3465.41
$ws.Range("F188").Value2 = 66
$ws.Range("G188").Value2 = 8801.76
$ws.Range("F238").Value2 = 23
$ws.Range("G238").Value2 = 1417.95
$ws.Range("F269").Value2 = 3
$ws.Range("G269").Value2 = 264.87
$ws.Range("F279").Value2 = 147
$ws.Range("G279").Value2 = 15502.62
$ws.Range("F281").Value2 = 120
$ws.Range("G281").Value2 = 8916
$ws.Range("F316").Value2 = 12
$ws.Range("G316").Value2 = 2756.16
$ws.Range("F339").Value2 = 209
$ws.Range("G339").Value2 = 8878.32
$ws.Range("F358").Value2 = 35
$ws.Range("G358").Value2 = 3398.5
$ws.Range("F368").Value2 = 22
$ws.Range("G368").Value2 = 2226.4
$ws.Range("F459").Value2 = 396
$ws.Range("G459").Value2 = 55673.64
$ws.Range("F471").Value2 = 340
$ws.Range("G471").Value2 = 56436.6
$ws.Range("F508").Value2 = 20
$ws.Range("G508").Value2 = 474
$ws.Range("F511").Value2 = 278
$ws.Range("G511").Value2 = 4406.3
$ws.Range("F517").Value2 = 143
$ws.Range("G517").Value2 = 8232.51
$ws.Range("F580").Value2 = 77
$ws.Range("G580").Value2 = 2070.53
$ws.Range("F582").Value2 = 259
$ws.Range("G582").Value2 = 2496.76
$ws.Range("F608").Value2 = 71
$ws.Range("G608").Value2 = 19256.62
$ws.Range("F692").Value2 = 209
$ws.Range("G692").Value2 = 17892.49
$ws.Range("F703").Value2 = 66
$ws.Range("G703").Value2 = 2849.88
$ws.Range("F720").Value2 = 86
$ws.Range("G720").Value2 = 5323.4
$ws.Range("F819").Value2 = 95
$ws.Range("G819").Value2 = 4571.4
$ws.Range("F825").Value2 = 448
$ws.Range("G825").Value2 = 35208.32
$ws.Range("F879").Value2 = 76
$ws.Range("G879").Value2 = 2718.52
$ws.Range("F882").Value2 = 18
$ws.Range("G882").Value2 = 803.16
$ws.Range("F891").Value2 = 1776
$ws.Range("G891").Value2 = 289683.36
$ws.Range("F893").Value2 = 87
$ws.Range("G893").Value2 = 24609.69
$ws.Range("F894").Value2 = 77
$ws.Range("G894").Value2 = 11138.05
$ws.Range("F909").Value2 = 34
$ws.Range("G909").Value2 = 5411.44
$ws.Range("F916").Value2 = 16
$ws.Range("G916").Value2 = 6476.32
$ws.Range("F933").Value2 = 20
$ws.Range("G933").Value2 = 4281.2

# --- Row pairs whose stock entries were swapped --------------------------
# Rows 151 <-> 152
$ws.Range("B151").Value2 = 64196
$ws.Range("F151").Value2 = 1
$ws.Range("G151").Value2 = 32143.58
$ws.Range("B152").Value2 = 65258
$ws.Range("F152").Value2 = 2
$ws.Range("G152").Value2 = 64287.16
# Rows 198 <-> 199
$ws.Range("B198").Value2 = 57756
$ws.Range("E198").Value2 = 79.37
$ws.Range("F198").Value2 = -100
$ws.Range("G198").Value2 = -6644
$ws.Range("B199").Value2 = 64350
$ws.Range("E199").Value2 = 70.63
$ws.Range("F199").Value2 = 2
$ws.Range("G199").Value2 = 132.88
# Rows 372 <-> 373
$ws.Range("B372").Value2 = 64985
$ws.Range("C372").Value2 = 'HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S'
$ws.Range("F372").Value2 = 13
$ws.Range("G372").Value2 = 1140.1
$ws.Range("B373").Value2 = 66196
$ws.Range("C373").Value2 = 'HIM-Total Care Baby Pants Drapers-Xl-9S'
$ws.Range("F373").Value2 = 29
$ws.Range("G373").Value2 = 2543.3
# Rows 567 <-> 568
$ws.Range("B567").Value2 = 64925
$ws.Range("E567").Value2 = 13.97
$ws.Range("F567").Value2 = 111
$ws.Range("G567").Value2 = 1459.65
$ws.Range("B568").Value2 = 45709
$ws.Range("E568").Value2 = 15.69
$ws.Range("F568").Value2 = -300
$ws.Range("G568").Value2 = -3945
# Rows 572 <-> 573
$ws.Range("B572").Value2 = 53595
$ws.Range("E572").Value2 = 17.61
$ws.Range("F572").Value2 = -335
$ws.Range("G572").Value2 = -4934.55
$ws.Range("B573").Value2 = 65067
$ws.Range("E573").Value2 = 15.65
$ws.Range("F573").Value2 = 126
$ws.Range("G573").Value2 = 1855.98
# Rows 672 <-> 673
$ws.Range("B672").Value2 = 60022
$ws.Range("E672").Value2 = 37.22
$ws.Range("F672").Value2 = -113
$ws.Range("G672").Value2 = -3709.79
$ws.Range("B673").Value2 = 64830
$ws.Range("E673").Value2 = 34.9
$ws.Range("F673").Value2 = 91
$ws.Range("G673").Value2 = 2987.53

# --- Sub Total / Grand Total rows, refreshed to match the corrections ---
$ws.Range("B41").Value2 = 81115.22
$ws.Range("B74").Value2 = 283760.23
$ws.Range("B102").Value2 = 132575.1
$ws.Range("B147").Value2 = 96435.00999999999
$ws.Range("B192").Value2 = 48985.2
$ws.Range("B250").Value2 = 100363.7
$ws.Range("B273").Value2 = 9289.719999999999
$ws.Range("B283").Value2 = 112846.25
$ws.Range("B318").Value2 = 23094.15
$ws.Range("B375").Value2 = 174655.22
$ws.Range("B461").Value2 = 121815.87
$ws.Range("B473").Value2 = 99080.92999999999
$ws.Range("B522").Value2 = 206037.55
$ws.Range("B587").Value2 = 41481.02
$ws.Range("B615").Value2 = 146272.39
$ws.Range("B695").Value2 = 185757.74
$ws.Range("B708").Value2 = 40701.31
$ws.Range("B732").Value2 = 38175.22
$ws.Range("B838").Value2 = 323799.41
$ws.Range("B885").Value2 = 25273.63
$ws.Range("B897").Value2 = 332687.91
$ws.Range("B913").Value2 = 19832.5
$ws.Range("B936").Value2 = 114502.42
$ws.Range("B942").Value2 = 5027978.76
$ws.Range("B943").Value2 = 5027978.76

